$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - (Intercept)
$ws.Range("B2").Value = 3241.773955
$ws.Range("D2").Value = 63.004168

# Row 3 - household_group_collapsed
$ws.Range("B3").Value = 1874.521102
$ws.Range("D3").Value = 18.215743
$ws.Range("E3").Value = 0

# Row 4 - Residuals
$ws.Range("B4").Value = 17082.504045
$ws.Range("C4").Value = 332

# Row 5 - SM-Control
$ws.Range("G5").Value = -2.914123
$ws.Range("H5").Value = -5.514826
$ws.Range("I5").Value = -0.31342
$ws.Range("J5").Value = 0.023669

# Row 6 - SM + Traps-Control
$ws.Range("G6").Value = 2.316188
$ws.Range("H6").Value = -0.415432
$ws.Range("I6").Value = 5.047808
$ws.Range("J6").Value = 0.114706

# Row 7 - SM + Traps-SM
$ws.Range("G7").Value = 5.230311
$ws.Range("H7").Value = 3.176315
$ws.Range("I7").Value = 7.284307
$ws.Range("J7").Value = 0
